$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates
$ws.Range("I2").Value = 5
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 7
$ws.Range("X2").Value = 7
$ws.Range("AB2").Value = 41
$ws.Range("AZ2").Value = 126

# Row 3 updates
$ws.Range("O3").Value = 1.33
$ws.Range("P3").Value = 3.25
